$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.4661921708185053
$ws1.Range("C2").Value = 0.075
$ws1.Range("D2").Value = 0.8571428571428571
$ws1.Range("E2").Value = 0.1379310344827586
$ws1.Range("F2").Value = 0.2777777777777778
$ws1.Range("G2").Value = 0.611764705882353
$ws1.Range("H2").Value = 0.6387774210807919
$ws1.Range("I2").Value = 24
$ws1.Range("J2").Value = 296
$ws1.Range("K2").Value = 238
$ws1.Range("L2").Value = 4

# --- Sheet: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.9834710743801653
$ws2.Range("C2").Value = 0.4456928838951311
$ws2.Range("D2").Value = 0.6134020618556701

$ws2.Range("B3").Value = 0.075
$ws2.Range("C3").Value = 0.8571428571428571
$ws2.Range("D3").Value = 0.1379310344827586

$ws2.Range("B4").Value = 0.4661921708185053
$ws2.Range("C4").Value = 0.4661921708185053
$ws2.Range("D4").Value = 0.4661921708185053
$ws2.Range("E4").Value = 0.4661921708185053

$ws2.Range("B5").Value = 0.5292355371900826
$ws2.Range("C5").Value = 0.6514178705189941
$ws2.Range("D5").Value = 0.3756665481692144

$ws2.Range("B6").Value = 0.9382091703185201
$ws2.Range("C6").Value = 0.4661921708185053
$ws2.Range("D6").Value = 0.5897131138726781

# --- Sheet: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 238
$ws3.Range("C2").Value = 296

$ws3.Range("B3").Value = 4
$ws3.Range("C3").Value = 24
